$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 33742.426
$ws.Range("I64").Value = 3185.7144
$ws.Range("J64").Value = 56257.895
$ws.Range("K64").Value = 3185.7144
$ws.Range("L64").Value = 56257.895
$ws.Range("M64").Value = -2937.7144
$ws.Range("N64").Value = -56753.895
$ws.Range("H67").Value = 33742.426
$ws.Range("I67").Value = 3185.7144
$ws.Range("J67").Value = 56257.895
$ws.Range("K67").Value = 3185.7144
$ws.Range("L67").Value = 56257.895
$ws.Range("M67").Value = -2327.7144
$ws.Range("N67").Value = -57973.895
$ws.Range("H87").Value = 13250
$ws.Range("J87").Value = 13250
$ws.Range("L87").Value = 13250
$ws.Range("N87").Value = -15746
$ws.Range("H90").Value = 13250
$ws.Range("J90").Value = 13250
$ws.Range("L90").Value = 39750
$ws.Range("N90").Value = -52230
$ws.Range("H111").Value = 578.5714
$ws.Range("J111").Value = 833
$ws.Range("L111").Value = 2499
$ws.Range("N111").Value = -8633
$ws.Range("H115").Value = 485.625
$ws.Range("I115").Value = 483.57144
$ws.Range("J115").Value = 500
$ws.Range("K115").Value = 1450.71432
$ws.Range("L115").Value = 1500
$ws.Range("M115").Value = 116.28568
$ws.Range("N115").Value = -4634
$ws.Range("H116").Value = 2057.946
$ws.Range("I116").Value = 1975.6129
$ws.Range("K116").Value = 1975.6129
$ws.Range("M116").Value = 1466.3871
$ws.Range("H118").Value = 523.1818
$ws.Range("I118").Value = 322.94116
$ws.Range("J118").Value = 1204
$ws.Range("K118").Value = 968.82348
$ws.Range("L118").Value = 3612
$ws.Range("M118").Value = 688.17652
$ws.Range("N118").Value = -6926
$ws.Range("H137").Value = 4006.7
$ws.Range("I137").Value = 4658.5625
$ws.Range("J137").Value = 1399.25
$ws.Range("K137").Value = 13975.6875
$ws.Range("L137").Value = 4197.75
$ws.Range("M137").Value = -11425.6875
$ws.Range("N137").Value = -9297.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2340.3333
$ws.Range("I2").Value = 1232.4
$ws.Range("J2").Value = 7880
$ws.Range("K2").Value = 1232.4
$ws.Range("L2").Value = 7880
$ws.Range("M2").Value = -1119.4
$ws.Range("N2").Value = -8106
$ws.Range("H88").Value = 2531.5386
$ws.Range("I88").Value = 2216.6667
$ws.Range("J88").Value = 2801.4285
$ws.Range("K88").Value = 2216.6667
$ws.Range("L88").Value = 2801.4285
$ws.Range("M88").Value = -1810.6667
$ws.Range("N88").Value = -3613.4285
$ws.Range("H91").Value = 2531.5386
$ws.Range("I91").Value = 2216.6667
$ws.Range("J91").Value = 2801.4285
$ws.Range("K91").Value = 2216.6667
$ws.Range("L91").Value = 2801.4285
$ws.Range("M91").Value = -812.6667000000002
$ws.Range("N91").Value = -5609.4285
$ws.Range("H116").Value = 2340.3333
$ws.Range("I116").Value = 1232.4
$ws.Range("J116").Value = 7880
$ws.Range("K116").Value = 1232.4
$ws.Range("L116").Value = 7880
$ws.Range("M116").Value = 1061.6
$ws.Range("N116").Value = -12468
$ws.Range("H132").Value = 21288.316
$ws.Range("I132").Value = 21801.12
$ws.Range("J132").Value = 17625.428
$ws.Range("K132").Value = 65403.36
$ws.Range("L132").Value = 52876.284
$ws.Range("M132").Value = -62873.36
$ws.Range("N132").Value = -57936.284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2340.3333
$ws.Range("I3").Value = 1232.4
$ws.Range("J3").Value = 7880
$ws.Range("K3").Value = 1232.4
$ws.Range("L3").Value = 7880
$ws.Range("M3").Value = -1118.4
$ws.Range("N3").Value = -8108
$ws.Range("H20").Value = 919.561
$ws.Range("I20").Value = 982.4091
$ws.Range("J20").Value = 846.7895
$ws.Range("K20").Value = 982.4091
$ws.Range("L20").Value = 846.7895
$ws.Range("M20").Value = -735.4091
$ws.Range("N20").Value = -1340.7895
$ws.Range("H86").Value = 31335.264
$ws.Range("I86").Value = 34726.324
$ws.Range("J86").Value = 16317.714
$ws.Range("K86").Value = 34726.324
$ws.Range("L86").Value = 16317.714
$ws.Range("M86").Value = -33603.324
$ws.Range("N86").Value = -18563.714
$ws.Range("H89").Value = 31335.264
$ws.Range("I89").Value = 34726.324
$ws.Range("J89").Value = 16317.714
$ws.Range("K89").Value = 173631.62
$ws.Range("L89").Value = 81588.57000000001
$ws.Range("M89").Value = -168015.62
$ws.Range("N89").Value = -92820.57000000001
$ws.Range("H107").Value = 1498.75
$ws.Range("I107").Value = 1331.6666
$ws.Range("K107").Value = 1331.6666
$ws.Range("M107").Value = 588.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 64496.906
$ws.Range("I31").Value = 73801.94500000001
$ws.Range("J31").Value = 8666.666999999999
$ws.Range("K31").Value = 73801.94500000001
$ws.Range("L31").Value = 8666.666999999999
$ws.Range("M31").Value = -73506.94500000001
$ws.Range("N31").Value = -9256.666999999999
$ws.Range("H34").Value = 64496.906
$ws.Range("I34").Value = 73801.94500000001
$ws.Range("J34").Value = 8666.666999999999
$ws.Range("K34").Value = 73801.94500000001
$ws.Range("L34").Value = 8666.666999999999
$ws.Range("M34").Value = -73599.94500000001
$ws.Range("N34").Value = -9070.666999999999
$ws.Range("H62").Value = 5464.625
$ws.Range("I62").Value = 4141
$ws.Range("J62").Value = 7670.6665
$ws.Range("K62").Value = 4141
$ws.Range("L62").Value = 7670.6665
$ws.Range("M62").Value = -3517
$ws.Range("N62").Value = -8918.666499999999
$ws.Range("H63").Value = 50271
$ws.Range("J63").Value = 50271
$ws.Range("L63").Value = 50271
$ws.Range("N63").Value = -51643
$ws.Range("H65").Value = 5464.625
$ws.Range("I65").Value = 4141
$ws.Range("J65").Value = 7670.6665
$ws.Range("K65").Value = 20705
$ws.Range("L65").Value = 38353.3325
$ws.Range("M65").Value = -17585
$ws.Range("N65").Value = -44593.3325
$ws.Range("H66").Value = 50271
$ws.Range("J66").Value = 50271
$ws.Range("L66").Value = 150813
$ws.Range("N66").Value = -157677
$ws.Range("H134").Value = 8879.678
$ws.Range("I134").Value = 5509.115
$ws.Range("J134").Value = 26406.6
$ws.Range("K134").Value = 16527.345
$ws.Range("L134").Value = 79219.79999999999
$ws.Range("M134").Value = -13992.345
$ws.Range("N134").Value = -84289.79999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 17816.455
$ws.Range("I122").Value = 622.75
$ws.Range("K122").Value = 5604.75
$ws.Range("M122").Value = -3154.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 61265.06
$ws.Range("I70").Value = 45960.5
$ws.Range("J70").Value = 94656.82000000001
$ws.Range("K70").Value = 45960.5
$ws.Range("L70").Value = 94656.82000000001
$ws.Range("M70").Value = -45690.5
$ws.Range("N70").Value = -95196.82000000001
$ws.Range("H73").Value = 61265.06
$ws.Range("I73").Value = 45960.5
$ws.Range("J73").Value = 94656.82000000001
$ws.Range("K73").Value = 45960.5
$ws.Range("L73").Value = 94656.82000000001
$ws.Range("M73").Value = -45024.5
$ws.Range("N73").Value = -96528.82000000001
$ws.Range("H97").Value = 47070.41
$ws.Range("I97").Value = 49216.094
$ws.Range("J97").Value = 2011
$ws.Range("K97").Value = 49216.094
$ws.Range("L97").Value = 2011
$ws.Range("M97").Value = -48720.094
$ws.Range("N97").Value = -3003
$ws.Range("H113").Value = 987.1429000000001
$ws.Range("I113").Value = 987.1429000000001
$ws.Range("K113").Value = 987.1429000000001
$ws.Range("M113").Value = 1182.8571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3021
$ws.Range("J7").Value = 3652.5
$ws.Range("L7").Value = 3652.5
$ws.Range("N7").Value = -3876.5
$ws.Range("H40").Value = 2791.3572
$ws.Range("I40").Value = 1958.4
$ws.Range("J40").Value = 4873.75
$ws.Range("K40").Value = 1958.4
$ws.Range("L40").Value = 4873.75
$ws.Range("M40").Value = -1822.4
$ws.Range("N40").Value = -5145.75
$ws.Range("H61").Value = 1391.4166
$ws.Range("I61").Value = 794
$ws.Range("K61").Value = 794
$ws.Range("M61").Value = -592
$ws.Range("H113").Value = 1391.4166
$ws.Range("I113").Value = 794
$ws.Range("K113").Value = 794
$ws.Range("M113").Value = 1376
$ws.Range("H126").Value = 3021
$ws.Range("J126").Value = 3652.5
$ws.Range("L126").Value = 10957.5
$ws.Range("N126").Value = -15897.5
